$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("JIND_raw")
$ws2 = $wb.Worksheets.Item("JIND")

# --- Sheet "JIND_raw": update recomputed confusion-matrix values ---
$ws1.Range("C2").Value = 0.0004580852038479157
$ws1.Range("G2").Value = 0.003216911764705882
$ws1.Range("C3").Value = 0.9656436097114063
$ws1.Range("D3").Value = 0.028142589118198873
$ws1.Range("C4").Value = 0.03344021988089785
$ws1.Range("D4").Value = 0.9690431519699813
$ws1.Range("G4").Value = 0.0
$ws1.Range("I4").Value = 0.006896551724137931
$ws1.Range("F6").Value = 0.7017543859649122
$ws1.Range("C7").Value = 0.0004580852038479157
$ws1.Range("F7").Value = 0.2982456140350877
$ws1.Range("G7").Value = 0.9871323529411765
$ws1.Range("H7").Value = 0.08169014084507042
$ws1.Range("G8").Value = 0.006893382352941176
$ws1.Range("H8").Value = 0.9183098591549296
$ws1.Range("I9").Value = 0.993103448275862

# --- Sheet "JIND": insert a new row for "Plasmacytoid dendritic cell" at row 10 ---
# (shifts old row 10 "Unassigned" down to row 11) and refresh all values
$ws2.Rows(10).Insert()

$ws2.Range("A2").Value = "B cell"
$ws2.Range("B2").Value = 0.9675767918088737
$ws2.Range("C2").Value = 0.0
$ws2.Range("D2").Value = 0.0009380863039399625
$ws2.Range("E2").Value = 0.0
$ws2.Range("F2").Value = 0.0
$ws2.Range("G2").Value = 0.0
$ws2.Range("H2").Value = 0.0
$ws2.Range("I2").Value = 0.0
$ws2.Range("J2").Value = 0.0
$ws2.Range("A3").Value = "CD4 T cell"
$ws2.Range("B3").Value = 0.0
$ws2.Range("C3").Value = 0.945487860742098
$ws2.Range("D3").Value = 0.021575984990619138
$ws2.Range("E3").Value = 0.0
$ws2.Range("F3").Value = 0.0
$ws2.Range("G3").Value = 0.00045955882352941176
$ws2.Range("H3").Value = 0.0
$ws2.Range("I3").Value = 0.0
$ws2.Range("J3").Value = 0.0
$ws2.Range("A4").Value = "CD8 T cell"
$ws2.Range("B4").Value = 0.0
$ws2.Range("C4").Value = 0.02244617498854787
$ws2.Range("D4").Value = 0.9652908067542214
$ws2.Range("E4").Value = 0.0
$ws2.Range("F4").Value = 0.0
$ws2.Range("G4").Value = 0.0
$ws2.Range("H4").Value = 0.0
$ws2.Range("I4").Value = 0.0
$ws2.Range("J4").Value = 0.0
$ws2.Range("A5").Value = "Hematopoietic stem cell"
$ws2.Range("B5").Value = 0.0
$ws2.Range("C5").Value = 0.0
$ws2.Range("D5").Value = 0.0
$ws2.Range("E5").Value = 0.8571428571428571
$ws2.Range("F5").Value = 0.0
$ws2.Range("G5").Value = 0.0
$ws2.Range("H5").Value = 0.0
$ws2.Range("I5").Value = 0.0
$ws2.Range("J5").Value = 0.0
$ws2.Range("A6").Value = "Megakaryocyte"
$ws2.Range("B6").Value = 0.0
$ws2.Range("C6").Value = 0.0
$ws2.Range("D6").Value = 0.0
$ws2.Range("E6").Value = 0.0
$ws2.Range("F6").Value = 0.6491228070175439
$ws2.Range("G6").Value = 0.0
$ws2.Range("H6").Value = 0.0
$ws2.Range("I6").Value = 0.0
$ws2.Range("J6").Value = 0.0
$ws2.Range("A7").Value = "Monocyte_CD14"
$ws2.Range("B7").Value = 0.0008532423208191126
$ws2.Range("C7").Value = 0.0
$ws2.Range("D7").Value = 0.0
$ws2.Range("E7").Value = 0.0
$ws2.Range("F7").Value = 0.22807017543859648
$ws2.Range("G7").Value = 0.9659926470588235
$ws2.Range("H7").Value = 0.036619718309859155
$ws2.Range("I7").Value = 0.0
$ws2.Range("J7").Value = 0.0
$ws2.Range("A8").Value = "Monocyte_FCGR3A"
$ws2.Range("B8").Value = 0.0
$ws2.Range("C8").Value = 0.0
$ws2.Range("D8").Value = 0.0
$ws2.Range("E8").Value = 0.0
$ws2.Range("F8").Value = 0.0
$ws2.Range("G8").Value = 0.0009191176470588235
$ws2.Range("H8").Value = 0.856338028169014
$ws2.Range("I8").Value = 0.0
$ws2.Range("J8").Value = 0.0
$ws2.Range("A9").Value = "NK cell"
$ws2.Range("B9").Value = 0.0
$ws2.Range("C9").Value = 0.0
$ws2.Range("D9").Value = 0.0
$ws2.Range("E9").Value = 0.0
$ws2.Range("F9").Value = 0.0
$ws2.Range("G9").Value = 0.00045955882352941176
$ws2.Range("H9").Value = 0.0
$ws2.Range("I9").Value = 0.9655172413793104
$ws2.Range("J9").Value = 0.0
$ws2.Range("A10").Value = "Plasmacytoid dendritic cell"
$ws2.Range("B10").Value = 0.0
$ws2.Range("C10").Value = 0.0
$ws2.Range("D10").Value = 0.0
$ws2.Range("E10").Value = 0.0
$ws2.Range("F10").Value = 0.0
$ws2.Range("G10").Value = 0.0
$ws2.Range("H10").Value = 0.0
$ws2.Range("I10").Value = 0.0
$ws2.Range("J10").Value = 0.9444444444444444
$ws2.Range("A11").Value = "Unassigned"
$ws2.Range("B11").Value = 0.031569965870307165
$ws2.Range("C11").Value = 0.0320659642693541
$ws2.Range("D11").Value = 0.012195121951219513
$ws2.Range("E11").Value = 0.14285714285714285
$ws2.Range("F11").Value = 0.12280701754385964
$ws2.Range("G11").Value = 0.03216911764705882
$ws2.Range("H11").Value = 0.10704225352112676
$ws2.Range("I11").Value = 0.034482758620689655
$ws2.Range("J11").Value = 0.05555555555555555
